# "working on big chart" - add portfolio total / percentage summary rows
# below the existing Old house .. Grain Field table (rows 2-10), and widen
# columns A and H to make room for the new labels / upcoming chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: "Total of all" - column sums of the B1:G10 table -------------
$ws.Range("A11").Value = "Total of all"
$ws.Range("B11").Value = 6275000
$ws.Range("C11").Value = 6728000
$ws.Range("D11").Value = 7354000
$ws.Range("E11").Value = 7435500
$ws.Range("F11").Value = 7468000
$ws.Range("G11").Value = 8083000

# C11 picks up a distinct (Calibri 12 black) font run in the source edit
$ws.Range("C11").Font.Name = "Calibri"
$ws.Range("C11").Font.Size = 12
$ws.Range("C11").Font.Color = 0

# --- Rows 12-14: category subtotals (based on column G / year 2020) -------
$ws.Range("A12").Value = "Total of appartments"
$ws.Range("B12").Formula = "=SUM(G4,G5)"

$ws.Range("A13").Value = "Total of houses"
$ws.Range("B13").Formula = "=SUM(G2,G6,G7,G8,G9)"

$ws.Range("A14").Value = "Total of land"
$ws.Range("B14").Formula = "=SUM(G3,G10)"

# --- Rows 15-17: percentage of portfolio ----------------------------------
$ws.Range("A15").Value = "Land % in portfolio:"
$ws.Range("B15").Formula = "=SUM(B14*100/G11)"

$ws.Range("A16").Value = "Houses % in portfolio:"
$ws.Range("B16").Formula = "=SUM(B13*100/G11)"

$ws.Range("A17").Value = "Apartments % in portfolio:"
$ws.Range("B17").Formula = "=SUM(B12*100/G11)"

# --- Column widths: widen the label column, and reserve an extra column ---
# for the upcoming chart ----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.46
$ws.Columns.Item(8).ColumnWidth = 18.76

# --- Selection follows the newly entered total row -------------------------
$null = $ws.Range("B11").Select()
